# Doku für Iteration 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# G3: was a numeric value (2.3), now a text note
$ws.Range("G3").Value = "2/3 nur teilweise"

# H3: clear the cell (limit value removed for this iteration)
$ws.Range("H3").ClearContents()

# Selection moves to C4
$ws.Range("C4").Select()
